# Adds two new columns ("Operacion" and "Zona") to the "NEW" sheet,
# matching the commit's automatic index.html/Excel export update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, same header text style (bold, bordered, centered) as
# the existing header row (A1:N1) -- copy formats only from A1 so we reuse
# the same style record instead of creating a near-duplicate one.
$ws.Range("O1").Value = "Operacion"
$ws.Range("P1").Value = "Zona"
$ws.Range("A1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Per-row "Operacion" (neighbourhood) / "Zona" values.
$data = @(
    @(2, "Saavedra", "Capital Norte"),
    @(3, "Paternal", "Capital Norte"),
    @(4, "Saavedra", "Capital Norte"),
    @(5, "Saavedra", "Capital Norte"),
    @(6, "Paternal", "Capital Norte"),
    @(7, "Palermo", "Capital Sur"),
    @(8, "Colegiales", "Capital Norte"),
    @(9, "Saavedra", "Capital Norte"),
    @(10, "Paternal", "Capital Norte"),
    @(11, "Colegiales", "Capital Norte"),
    @(12, "Colegiales", "Capital Norte"),
    @(13, "Saavedra", "Capital Norte"),
    @(14, "Colegiales", "Capital Norte"),
    @(15, "Paternal", "Capital Norte"),
    @(16, "Paternal", "Capital Norte"),
    @(17, "Saavedra", "Capital Norte"),
    @(18, "Saavedra", "Capital Norte"),
    @(19, "Colegiales", "Capital Norte"),
    @(20, "Saavedra", "Capital Norte"),
    @(21, "Paternal", "Capital Norte"),
    @(22, "Colegiales", "Capital Norte"),
    @(23, "Paternal", "Capital Norte"),
    @(24, "Colegiales", "Capital Norte"),
    @(25, "Palermo", "Capital Sur"),
    @(26, "Paternal", "Capital Norte"),
    @(27, "Paternal", "Capital Norte"),
    @(28, "Paternal", "Capital Norte"),
    @(29, "Paternal", "Capital Norte"),
    @(30, "Palermo", "Capital Sur"),
    @(31, "Saavedra", "Capital Norte"),
    @(32, "Paternal", "Capital Norte"),
    @(33, "Paternal", "Capital Norte"),
    @(34, "Saavedra", "Capital Norte"),
    @(35, "Colegiales", "Capital Norte"),
    @(36, "Saavedra", "Capital Norte"),
    @(37, "Palermo", "Capital Sur"),
    @(38, "Paternal", "Capital Norte"),
    @(39, "Paternal", "Capital Norte"),
    @(40, "Paternal", "Capital Norte"),
    @(41, "Paternal", "Capital Norte"),
    @(42, "Palermo", "Capital Sur"),
    @(43, "Palermo", "Capital Sur"),
    @(44, "Palermo", "Capital Sur"),
    @(45, "Palermo", "Capital Sur"),
    @(46, "Saavedra", "Capital Norte"),
    @(47, "Saavedra", "Capital Norte"),
    @(48, "Palermo", "Capital Sur"),
    @(49, "Saavedra", "Capital Norte")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 15).Value = $row[1]
    $ws.Cells.Item($r, 16).Value = $row[2]
}
